# Insert a new data row into the worksheet between the current row 54 and
# row 55, shifting all subsequent rows down by one (rows 55-131 become
# rows 56-132), and fill the newly inserted row 55 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 55 (pushes existing row 55.. down by 1)
$ws.Rows.Item(55).Insert()

# Populate the new row 55 with its values
$ws.Cells.Item(55, 1).Value2 = 4
$ws.Cells.Item(55, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(55, 3).Value2 = "Los Lagos"
$ws.Cells.Item(55, 4).Value2 = 44792
$ws.Cells.Item(55, 5).Value2 = 10
$ws.Cells.Item(55, 6).Value2 = 100112022
$ws.Cells.Item(55, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(55, 8).Value2 = "Perfection"
$ws.Cells.Item(55, 9).Value2 = "Primera"
$ws.Cells.Item(55, 10).Value2 = 40
$ws.Cells.Item(55, 11).Value2 = 46000
$ws.Cells.Item(55, 12).Value2 = 46000
$ws.Cells.Item(55, 13).Value2 = 46000
$ws.Cells.Item(55, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(55, 15).Value2 = "Provincia de Huasco"
$ws.Cells.Item(55, 16).Value2 = 1840
$ws.Cells.Item(55, 17).Value2 = 25
$ws.Cells.Item(55, 18).Value2 = "Hortaliza"
